$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill AR column for rows 156-158 ---
$ws.Range("AR156").Value = 17.2896458
$ws.Range("AR157").Value = 10.6618117
$ws.Range("AR158").Value = 13.1818351

# --- Corrected values in row 178 (re-estimated after adding new survey days) ---
$ws.Range("M178").Value = 33.941274
$ws.Range("W178").Value = 11.8515132
$ws.Range("AO178").Value = 30.4751199
$ws.Range("AV178").Value = 30.1001694
$ws.Range("AW178").Value = 34.965013
$ws.Range("BB178").Value = 17.3627136

# --- Row 179: fill in data for "27 07 2020" (date label cell A179 already set) ---
$ws.Range("B179").Value = 23.2722731
$ws.Range("C179").Value = 37.207655
$ws.Range("D179").Value = 28.1990124
$ws.Range("F179").Value = 31.3528597
$ws.Range("G179").Value = 21.3723123
$ws.Range("H179").Value = 18.2413855
$ws.Range("I179").Value = 10.5285665
$ws.Range("J179").Value = 16.5088757
$ws.Range("K179").Value = 15.8197217
$ws.Range("L179").Value = 29.5329248
$ws.Range("M179").Value = 34.3384716
$ws.Range("O179").Value = 13.1493506
$ws.Range("P179").Value = 28.9971283
$ws.Range("Q179").Value = 33.1745067
$ws.Range("R179").Value = 18.1922995
$ws.Range("S179").Value = 23.0880374
$ws.Range("T179").Value = 27.7967528
$ws.Range("U179").Value = 24.7203455
$ws.Range("V179").Value = 36.3187442
$ws.Range("W179").Value = 11.5554747
$ws.Range("X179").Value = 16.1430649
$ws.Range("Y179").Value = 10.4112694
$ws.Range("Z179").Value = 16.117721
$ws.Range("AA179").Value = 20.7499974
$ws.Range("AB179").Value = 25.8491547
$ws.Range("AD179").Value = 38.0403174
$ws.Range("AE179").Value = 24.1585283
$ws.Range("AF179").Value = 21.6722813
$ws.Range("AG179").Value = 26.1865129
$ws.Range("AH179").Value = 25.5745397
$ws.Range("AI179").Value = 10.3998597
$ws.Range("AJ179").Value = 12.4438023
$ws.Range("AK179").Value = 19.5587566
$ws.Range("AL179").Value = 24.7343903
$ws.Range("AM179").Value = 11.9620868
$ws.Range("AN179").Value = 21.2973936
$ws.Range("AO179").Value = 29.9100179
$ws.Range("AP179").Value = 16.5792233
$ws.Range("AQ179").Value = 15.115563
$ws.Range("AS179").Value = 13.4042867
$ws.Range("AT179").Value = 32.1506334
$ws.Range("AU179").Value = 23.1461725
$ws.Range("AV179").Value = 30.2285282
$ws.Range("AW179").Value = 34.3201445
$ws.Range("AX179").Value = 26.5992246
$ws.Range("AY179").Value = 19.3270404
$ws.Range("BA179").Value = 8.157685
$ws.Range("BB179").Value = 17.1632997
$ws.Range("BC179").Value = 21.2253494
$ws.Range("BD179").Value = 21.9692679
$ws.Range("BE179").Value = 17.3675279

# --- Row 180: "28 07 2020" ---
$ws.Range("A180").Value = "28 07 2020"
$ws.Range("B180").Value = 23.3535529
$ws.Range("C180").Value = 37.5869583
$ws.Range("D180").Value = 28.9426303
$ws.Range("F180").Value = 30.840195
$ws.Range("G180").Value = 21.3034525
$ws.Range("H180").Value = 18.1961511
$ws.Range("I180").Value = 10.1993068
$ws.Range("J180").Value = 16.3647343
$ws.Range("K180").Value = 16.2661196
$ws.Range("L180").Value = 29.0400876
$ws.Range("M180").Value = 33.9439115
$ws.Range("O180").Value = 13.0089059
$ws.Range("P180").Value = 28.3865655
$ws.Range("Q180").Value = 33.0228138
$ws.Range("R180").Value = 18.1438213
$ws.Range("S180").Value = 23.015616
$ws.Range("T180").Value = 27.6235108
$ws.Range("U180").Value = 24.6648119
$ws.Range("V180").Value = 35.8821955
$ws.Range("W180").Value = 11.3176299
$ws.Range("X180").Value = 16.6176227
$ws.Range("Y180").Value = 11.2726651
$ws.Range("Z180").Value = 16.3361849
$ws.Range("AA180").Value = 21.1270903
$ws.Range("AB180").Value = 26.0255015
$ws.Range("AD180").Value = 36.9573422
$ws.Range("AE180").Value = 25.3901903
$ws.Range("AF180").Value = 21.7194777
$ws.Range("AG180").Value = 26.1150565
$ws.Range("AH180").Value = 26.0925525
$ws.Range("AI180").Value = 10.7970506
$ws.Range("AJ180").Value = 12.2027682
$ws.Range("AK180").Value = 19.4235261
$ws.Range("AL180").Value = 24.671335
$ws.Range("AM180").Value = 11.9630724
$ws.Range("AN180").Value = 21.2541906
$ws.Range("AO180").Value = 30.037442
$ws.Range("AP180").Value = 16.805486
$ws.Range("AQ180").Value = 14.9360171
$ws.Range("AS180").Value = 13.1231464
$ws.Range("AT180").Value = 31.4405498
$ws.Range("AU180").Value = 21.7611024
$ws.Range("AV180").Value = 30.189567
$ws.Range("AW180").Value = 34.0038966
$ws.Range("AX180").Value = 26.134996
$ws.Range("AY180").Value = 19.2295393
$ws.Range("BA180").Value = 8.5203696
$ws.Range("BB180").Value = 17.5279079
$ws.Range("BC180").Value = 21.4613191
$ws.Range("BD180").Value = 21.0156045
$ws.Range("BE180").Value = 17.4617931

# --- Row 181: "29 07 2020" ---
$ws.Range("A181").Value = "29 07 2020"
$ws.Range("B181").Value = 22.6538462
$ws.Range("C181").Value = 38.0744638
$ws.Range("D181").Value = 28.9353308
$ws.Range("F181").Value = 30.2467095
$ws.Range("G181").Value = 21.189032
$ws.Range("H181").Value = 18.2528109
$ws.Range("I181").Value = 10.1231783
$ws.Range("J181").Value = 15.1678657
$ws.Range("K181").Value = 16.6369394
$ws.Range("L181").Value = 28.7196227
$ws.Range("M181").Value = 33.7977871
$ws.Range("O181").Value = 13.0363036
$ws.Range("P181").Value = 27.5866433
$ws.Range("Q181").Value = 33.0303248
$ws.Range("R181").Value = 18.2868178
$ws.Range("S181").Value = 23.0304971
$ws.Range("T181").Value = 26.9522131
$ws.Range("U181").Value = 24.4379686
$ws.Range("V181").Value = 35.6824533
$ws.Range("W181").Value = 11.5412586
$ws.Range("X181").Value = 16.2966909
$ws.Range("Y181").Value = 10.3165651
$ws.Range("Z181").Value = 16.3947097
$ws.Range("AA181").Value = 21.0320772
$ws.Range("AB181").Value = 26.0832501
$ws.Range("AD181").Value = 36.7164578
$ws.Range("AE181").Value = 24.1969618
$ws.Range("AF181").Value = 21.8057148
$ws.Range("AG181").Value = 25.2861559
$ws.Range("AH181").Value = 26.4557452
$ws.Range("AI181").Value = 10.8064986
$ws.Range("AJ181").Value = 12.3002424
$ws.Range("AK181").Value = 19.7415905
$ws.Range("AL181").Value = 24.8328892
$ws.Range("AM181").Value = 11.8285604
$ws.Range("AN181").Value = 21.0828784
$ws.Range("AO181").Value = 29.5949861
$ws.Range("AP181").Value = 16.4016042
$ws.Range("AQ181").Value = 14.6768578
$ws.Range("AS181").Value = 13.7097463
$ws.Range("AT181").Value = 31.1087288
$ws.Range("AU181").Value = 22.7306608
$ws.Range("AV181").Value = 30.1394756
$ws.Range("AW181").Value = 33.5128643
$ws.Range("AX181").Value = 26.6742385
$ws.Range("AY181").Value = 19.031965
$ws.Range("BA181").Value = 7.913961
$ws.Range("BB181").Value = 17.6198862
$ws.Range("BC181").Value = 21.609292
$ws.Range("BD181").Value = 22.0214023
$ws.Range("BE181").Value = 19.0949462

# --- Row 182: "30 07 2020" ---
$ws.Range("A182").Value = "30 07 2020"
$ws.Range("B182").Value = 22.1991701
$ws.Range("C182").Value = 36.5938748
$ws.Range("D182").Value = 28.9641188
$ws.Range("F182").Value = 29.6811328
$ws.Range("G182").Value = 21.1364448
$ws.Range("H182").Value = 18.3577128
$ws.Range("I182").Value = 9.9527311
$ws.Range("J182").Value = 13.5118306
$ws.Range("K182").Value = 16.7149479
$ws.Range("L182").Value = 28.6191258
$ws.Range("M182").Value = 33.5149611
$ws.Range("O182").Value = 12.273033
$ws.Range("P182").Value = 27.9940068
$ws.Range("Q182").Value = 33.338876
$ws.Range("R182").Value = 18.2910731
$ws.Range("S182").Value = 23.4382866
$ws.Range("T182").Value = 26.6107881
$ws.Range("U182").Value = 24.901015
$ws.Range("V182").Value = 36.0120962
$ws.Range("W182").Value = 11.457054
$ws.Range("X182").Value = 16.4098504
$ws.Range("Y182").Value = 10.8931241
$ws.Range("Z182").Value = 16.4516333
$ws.Range("AA182").Value = 21.0501156
$ws.Range("AB182").Value = 25.6821463
$ws.Range("AD182").Value = 36.5499902
$ws.Range("AE182").Value = 25.26119
$ws.Range("AF182").Value = 21.5392802
$ws.Range("AG182").Value = 25.1260257
$ws.Range("AH182").Value = 26.2642518
$ws.Range("AI182").Value = 10.6395578
$ws.Range("AJ182").Value = 12.2963799
$ws.Range("AK182").Value = 19.2829876
$ws.Range("AL182").Value = 24.1256961
$ws.Range("AM182").Value = 11.7403202
$ws.Range("AN182").Value = 21.0345974
$ws.Range("AO182").Value = 29.2196056
$ws.Range("AP182").Value = 16.4370176
$ws.Range("AQ182").Value = 14.5271234
$ws.Range("AS182").Value = 13.4429931
$ws.Range("AT182").Value = 30.0801317
$ws.Range("AU182").Value = 23.5865784
$ws.Range("AV182").Value = 30.6738652
$ws.Range("AW182").Value = 33.119137
$ws.Range("AX182").Value = 26.9256295
$ws.Range("AY182").Value = 18.817062
$ws.Range("BA182").Value = 8.6054762
$ws.Range("BB182").Value = 17.8040298
$ws.Range("BC182").Value = 21.6342954
$ws.Range("BD182").Value = 22.3370754
$ws.Range("BE182").Value = 19.5756122

# --- Row 183: "31 07 2020" ---
$ws.Range("A183").Value = "31 07 2020"
$ws.Range("B183").Value = 23.6996644
$ws.Range("C183").Value = 36.882389
$ws.Range("D183").Value = 28.4247365
$ws.Range("F183").Value = 28.9708673
$ws.Range("G183").Value = 20.9905935
$ws.Range("H183").Value = 18.4982598
$ws.Range("I183").Value = 9.8568751
$ws.Range("J183").Value = 13.2653061
$ws.Range("K183").Value = 16.8789809
$ws.Range("L183").Value = 28.136063
$ws.Range("M183").Value = 33.0112949
$ws.Range("O183").Value = 13.2343651
$ws.Range("P183").Value = 28.8718424
$ws.Range("Q183").Value = 33.3502634
$ws.Range("R183").Value = 18.5549258
$ws.Range("S183").Value = 24.2456447
$ws.Range("T183").Value = 27.3062209
$ws.Range("U183").Value = 24.683321
$ws.Range("V183").Value = 36.0861423
$ws.Range("W183").Value = 11.4769717
$ws.Range("X183").Value = 16.3204095
$ws.Range("Y183").Value = 10.366511
$ws.Range("Z183").Value = 16.7361767
$ws.Range("AA183").Value = 20.7366986
$ws.Range("AB183").Value = 25.6147101
$ws.Range("AD183").Value = 36.9825016
$ws.Range("AE183").Value = 24.0077048
$ws.Range("AF183").Value = 21.4794335
$ws.Range("AG183").Value = 26.1290907
$ws.Range("AH183").Value = 25.5467493
$ws.Range("AI183").Value = 11.0631104
$ws.Range("AJ183").Value = 12.2391221
$ws.Range("AK183").Value = 19.8519229
$ws.Range("AL183").Value = 24.4464509
$ws.Range("AM183").Value = 11.5745927
$ws.Range("AN183").Value = 21.3513204
$ws.Range("AO183").Value = 29.4102707
$ws.Range("AP183").Value = 16.336303
$ws.Range("AQ183").Value = 14.5604111
$ws.Range("AS183").Value = 14.1410496
$ws.Range("AT183").Value = 30.0262553
$ws.Range("AU183").Value = 22.3931185
$ws.Range("AV183").Value = 30.2534635
$ws.Range("AW183").Value = 32.5307177
$ws.Range("AX183").Value = 26.6186647
$ws.Range("AY183").Value = 18.8578985
$ws.Range("BA183").Value = 8.137343
$ws.Range("BB183").Value = 18.0277284
$ws.Range("BC183").Value = 21.6508977
$ws.Range("BD183").Value = 21.0499795
$ws.Range("BE183").Value = 20.2939183

# --- Row 184: "01 08 2020" ---
$ws.Range("A184").Value = "01 08 2020"
$ws.Range("B184").Value = 21.2204724
$ws.Range("C184").Value = 36.5460141
$ws.Range("D184").Value = 28.2646909
$ws.Range("F184").Value = 28.4173855
$ws.Range("G184").Value = 20.8859919
$ws.Range("H184").Value = 18.5661379
$ws.Range("I184").Value = 9.7588877
$ws.Range("J184").Value = 13.2552404
$ws.Range("K184").Value = 16.704481
$ws.Range("L184").Value = 27.9296046
$ws.Range("M184").Value = 32.7520169
$ws.Range("O184").Value = 12.8764597
$ws.Range("P184").Value = 28.1422174
$ws.Range("Q184").Value = 32.6668298
$ws.Range("R184").Value = 18.5199444
$ws.Range("S184").Value = 24.3418949
$ws.Range("T184").Value = 27.2134091
$ws.Range("U184").Value = 24.5835014
$ws.Range("V184").Value = 36.2889455
$ws.Range("W184").Value = 11.3178522
$ws.Range("X184").Value = 16.2219007
$ws.Range("Y184").Value = 10.2540501
$ws.Range("Z184").Value = 16.5892754
$ws.Range("AA184").Value = 20.8838102
$ws.Range("AB184").Value = 25.6044215
$ws.Range("AD184").Value = 36.7971856
$ws.Range("AE184").Value = 24.6052772
$ws.Range("AF184").Value = 21.3074439
$ws.Range("AG184").Value = 25.9072332
$ws.Range("AH184").Value = 25.424786
$ws.Range("AI184").Value = 10.7748651
$ws.Range("AJ184").Value = 12.1843797
$ws.Range("AK184").Value = 19.4123646
$ws.Range("AL184").Value = 23.9569418
$ws.Range("AM184").Value = 11.5081398
$ws.Range("AN184").Value = 21.2748349
$ws.Range("AO184").Value = 28.7830837
$ws.Range("AP184").Value = 16.0630317
$ws.Range("AQ184").Value = 14.6474637
$ws.Range("AS184").Value = 13.6279005
$ws.Range("AT184").Value = 30.2957086
$ws.Range("AU184").Value = 24.1167379
$ws.Range("AV184").Value = 30.2964955
$ws.Range("AW184").Value = 32.3372983
$ws.Range("AX184").Value = 26.8830308
$ws.Range("AY184").Value = 18.7681472
$ws.Range("BA184").Value = 6.9927241
$ws.Range("BB184").Value = 17.7678466
$ws.Range("BC184").Value = 21.1001513
$ws.Range("BD184").Value = 21.8926953
$ws.Range("BE184").Value = 20.6200192

# --- Row 185: "02 08 2020" ---
$ws.Range("A185").Value = "02 08 2020"
$ws.Range("B185").Value = 20.9347997
$ws.Range("C185").Value = 36.6477283
$ws.Range("D185").Value = 28.7608513
$ws.Range("F185").Value = 28.1847403
$ws.Range("G185").Value = 20.7713579
$ws.Range("H185").Value = 18.4094009
$ws.Range("I185").Value = 9.4290071
$ws.Range("J185").Value = 12.7068558
$ws.Range("K185").Value = 17.4721318
$ws.Range("L185").Value = 27.6205901
$ws.Range("M185").Value = 32.4350923
$ws.Range("O185").Value = 13.1578947
$ws.Range("P185").Value = 28.3752279
$ws.Range("Q185").Value = 31.8936125
$ws.Range("R185").Value = 18.9735377
$ws.Range("S185").Value = 24.6284958
$ws.Range("T185").Value = 26.14837
$ws.Range("U185").Value = 24.9738224
$ws.Range("V185").Value = 35.3658584
$ws.Range("W185").Value = 11.0353969
$ws.Range("X185").Value = 15.9726297
$ws.Range("Y185").Value = 10.2188493
$ws.Range("Z185").Value = 16.6346916
$ws.Range("AA185").Value = 20.6548986
$ws.Range("AB185").Value = 25.4929815
$ws.Range("AD185").Value = 37.0307178
$ws.Range("AE185").Value = 25.3123052
$ws.Range("AF185").Value = 21.2671814
$ws.Range("AG185").Value = 25.3886855
$ws.Range("AH185").Value = 25.6816512
$ws.Range("AI185").Value = 10.3773585
$ws.Range("AJ185").Value = 11.9818565
$ws.Range("AK185").Value = 18.4527316
$ws.Range("AL185").Value = 23.4478399
$ws.Range("AM185").Value = 11.7444509
$ws.Range("AN185").Value = 21.6656597
$ws.Range("AO185").Value = 29.132302
$ws.Range("AP185").Value = 15.4569926
$ws.Range("AQ185").Value = 14.515818
$ws.Range("AS185").Value = 13.504236
$ws.Range("AT185").Value = 29.7023328
$ws.Range("AU185").Value = 22.4352593
$ws.Range("AV185").Value = 30.3552956
$ws.Range("AW185").Value = 31.7044423
$ws.Range("AX185").Value = 27.049002
$ws.Range("AY185").Value = 18.5329562
$ws.Range("BA185").Value = 8.3680238
$ws.Range("BB185").Value = 17.8136131
$ws.Range("BC185").Value = 21.0273407
$ws.Range("BD185").Value = 22.6753
$ws.Range("BE185").Value = 22.7153665

# --- Row 186: "03 08 2020" ---
$ws.Range("A186").Value = "03 08 2020"
$ws.Range("B186").Value = 20.7021792
$ws.Range("C186").Value = 36.4498671
$ws.Range("D186").Value = 29.0991779
$ws.Range("F186").Value = 27.2690791
$ws.Range("G186").Value = 20.6404813
$ws.Range("H186").Value = 18.3837893
$ws.Range("I186").Value = 9.3948771
$ws.Range("J186").Value = 12.6682987
$ws.Range("K186").Value = 16.4766558
$ws.Range("L186").Value = 27.4637072
$ws.Range("M186").Value = 32.4534326
$ws.Range("O186").Value = 13.9813193
$ws.Range("P186").Value = 28.1873626
$ws.Range("Q186").Value = 32.362727
$ws.Range("R186").Value = 19.1834241
$ws.Range("S186").Value = 24.7071328
$ws.Range("T186").Value = 26.3381532
$ws.Range("U186").Value = 24.8531019
$ws.Range("V186").Value = 35.1411825
$ws.Range("W186").Value = 11.2826416
$ws.Range("X186").Value = 15.6429842
$ws.Range("Y186").Value = 10.3867214
$ws.Range("Z186").Value = 16.8416103
$ws.Range("AA186").Value = 21.1931739
$ws.Range("AB186").Value = 25.4863863
$ws.Range("AD186").Value = 36.6711853
$ws.Range("AE186").Value = 24.4715227
$ws.Range("AF186").Value = 21.0663477
$ws.Range("AG186").Value = 26.188426
$ws.Range("AH186").Value = 26.3388391
$ws.Range("AI186").Value = 10.1258581
$ws.Range("AJ186").Value = 12.0852196
$ws.Range("AK186").Value = 19.0366994
$ws.Range("AL186").Value = 23.9734002
$ws.Range("AM186").Value = 11.7139513
$ws.Range("AN186").Value = 21.6547488
$ws.Range("AO186").Value = 29.3228571
$ws.Range("AP186").Value = 15.1720903
$ws.Range("AQ186").Value = 14.6270193
$ws.Range("AS186").Value = 14.0551491
$ws.Range("AT186").Value = 29.7531213
$ws.Range("AU186").Value = 23.7572351
$ws.Range("AV186").Value = 29.7842703
$ws.Range("AW186").Value = 31.2136246
$ws.Range("AX186").Value = 26.9166807
$ws.Range("AY186").Value = 18.7451017
$ws.Range("BA186").Value = 7.6656442
$ws.Range("BB186").Value = 17.489207
$ws.Range("BC186").Value = 20.7549912
$ws.Range("BD186").Value = 22.5386432
$ws.Range("BE186").Value = 22.1575372

# --- Rows 187-188: date labels only, no survey data yet ---
$ws.Range("A187").Value = "04 08 2020"
$ws.Range("A188").Value = "05 08 2020"
